# CW3M_NSantiam.xlsx regression-testing workbook update
# - Adds a new calibration run row (row 13) to the "2010-18" sheet
#   ("Baseline 2010-18 C402" - added 174 cfs spring in the Little Nsantiam basin)
# - Moves the active sheet/tab selection from "2000-09 spinup" to "2010-18"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2010-18")

# New calibration run row
$ws.Cells.Item(13, 1).Value = "CW3M"
$ws.Cells.Item(13, 2).Value = "Baseline 2010-18 C402"
$ws.Cells.Item(13, 3).Value = "2010-18"

$ws.Cells.Item(13, 4).Value = 529.14105211111109
$ws.Cells.Item(13, 4).NumberFormat = "0.00"
$ws.Cells.Item(13, 4).Interior.Color = 65535

$ws.Cells.Item(13, 5).Value = 2094.2995878888887
$ws.Cells.Item(13, 5).NumberFormat = "0.00"

$ws.Cells.Item(13, 6).Value = 1.6230948888888888
$ws.Cells.Item(13, 6).NumberFormat = "0.00"
$ws.Cells.Item(13, 6).Interior.Color = 65535

$ws.Cells.Item(13, 7).Value = 332.2750817777777
$ws.Cells.Item(13, 7).NumberFormat = "0.00"
$ws.Cells.Item(13, 7).Interior.Color = 65535

$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 8).NumberFormat = "0.00"

$ws.Cells.Item(13, 9).Value = 7.3481075555555559
$ws.Cells.Item(13, 9).NumberFormat = "0.00"
$ws.Cells.Item(13, 9).Interior.Color = 65535

$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 10).NumberFormat = "0.00"

$ws.Cells.Item(13, 11).Value = 520.64759333333336
$ws.Cells.Item(13, 11).NumberFormat = "0.00"
$ws.Cells.Item(13, 11).Interior.Color = 65535

$ws.Cells.Item(13, 12).Value = 91.777595333333338
$ws.Cells.Item(13, 12).NumberFormat = "0.00"

$ws.Cells.Item(13, 13).Value = 1836.442098
$ws.Cells.Item(13, 13).NumberFormat = "0.00"
$ws.Cells.Item(13, 13).Interior.Color = 65535

$ws.Cells.Item(13, 14).Value = 515.86576666666667
$ws.Cells.Item(13, 14).NumberFormat = "0.00"
$ws.Cells.Item(13, 14).Interior.Color = 65535

$ws.Cells.Item(13, 15).Value = 3819.0314398888886
$ws.Cells.Item(13, 15).NumberFormat = "0"
$ws.Cells.Item(13, 15).Interior.Color = 255
$ws.Cells.Item(13, 15).Font.Bold = $true

$ws.Cells.Item(13, 16).Value = 2216.8192002222222
$ws.Cells.Item(13, 16).NumberFormat = "0"

$ws.Cells.Item(13, 17).Value = 0.0461291111111111137854657
$ws.Cells.Item(13, 17).NumberFormat = "0.00"
$ws.Cells.Item(13, 17).Interior.Color = 65535

$ws.Cells.Item(13, 18).Value = -0.0000102222222222222179009
$ws.Cells.Item(13, 18).NumberFormat = "0.000000"
$ws.Cells.Item(13, 18).Interior.Color = 65535

$ws.Cells.Item(13, 20).Value = "added 174 cfs spring in the Little Nsantiam basin"

# Move the active tab / selection to the "2010-18" sheet (was "2000-09 spinup")
$ws.Select()
$ws.Range("K13").Select()
